$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 429.9091
$ws.Range("J2").Value = 608
$ws.Range("L2").Value = 608
$ws.Range("N2").Value = -834
$ws.Range("H62").Value = 6259.9
$ws.Range("I62").Value = 4650
$ws.Range("J62").Value = 7333.1665
$ws.Range("K62").Value = 4650
$ws.Range("L62").Value = 7333.1665
$ws.Range("M62").Value = -4026
$ws.Range("N62").Value = -8581.166499999999
$ws.Range("H65").Value = 6259.9
$ws.Range("I65").Value = 4650
$ws.Range("J65").Value = 7333.1665
$ws.Range("K65").Value = 23250
$ws.Range("L65").Value = 36665.8325
$ws.Range("M65").Value = -20130
$ws.Range("N65").Value = -42905.8325
$ws.Range("H135").Value = 1057.6666
$ws.Range("I135").Value = 583.1667
$ws.Range("K135").Value = 5248.5003
$ws.Range("M135").Value = -2713.5003
$ws.Range("H140").Value = 91666.664
$ws.Range("J140").Value = 91666.664
$ws.Range("L140").Value = 91666.664
$ws.Range("N140").Value = -102026.664
$ws.Range("H141").Value = 2470.5
$ws.Range("I141").Value = 2114.6
$ws.Range("J141").Value = 4250
$ws.Range("K141").Value = 6343.799999999999
$ws.Range("L141").Value = 12750
$ws.Range("M141").Value = -1163.799999999999
$ws.Range("N141").Value = -23110

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 368.35715
$ws.Range("I2").Value = 375.53845
$ws.Range("K2").Value = 375.53845
$ws.Range("M2").Value = -262.53845
$ws.Range("H45").Value = 1870
$ws.Range("I45").Value = 1870
$ws.Range("K45").Value = 1870
$ws.Range("M45").Value = -1493
$ws.Range("H61").Value = 1377.1305
$ws.Range("I61").Value = 1377.1305
$ws.Range("K61").Value = 1377.1305
$ws.Range("M61").Value = -1165.1305
$ws.Range("H116").Value = 368.35715
$ws.Range("I116").Value = 375.53845
$ws.Range("K116").Value = 375.53845
$ws.Range("M116").Value = 1918.46155
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H136").Value = 1377.1305
$ws.Range("I136").Value = 1377.1305
$ws.Range("K136").Value = 4131.3915
$ws.Range("M136").Value = -1581.3915

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 368.35715
$ws.Range("I3").Value = 375.53845
$ws.Range("K3").Value = 375.53845
$ws.Range("M3").Value = -261.53845
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 80.25
$ws.Range("I7").Value = 102
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 102
$ws.Range("L7").Value = 15
$ws.Range("M7").Value = 11
$ws.Range("N7").Value = -241
$ws.Range("H122").Value = 5996.926
$ws.Range("I122").Value = 5915.0625
$ws.Range("J122").Value = 6116
$ws.Range("K122").Value = 17745.1875
$ws.Range("L122").Value = 18348
$ws.Range("M122").Value = -15295.1875
$ws.Range("N122").Value = -23248

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91103
$ws.Range("J2").Value = 204
$ws.Range("L2").Value = 1224
$ws.Range("N2").Value = -1450
$ws.Range("H4").Value = 133367896
$ws.Range("I4").Value = 181864000
$ws.Range("J4").Value = 3624.25
$ws.Range("K4").Value = 545592000
$ws.Range("L4").Value = 10872.75
$ws.Range("M4").Value = -545591888
$ws.Range("N4").Value = -11096.75
$ws.Range("H29").Value = 381.2
$ws.Range("I29").Value = 398
$ws.Range("J29").Value = 364.4
$ws.Range("K29").Value = 1194
$ws.Range("L29").Value = 1093.2
$ws.Range("M29").Value = -917
$ws.Range("N29").Value = -1647.2
$ws.Range("H34").Value = 1483.0769
$ws.Range("J34").Value = 2126.8572
$ws.Range("L34").Value = 6380.571599999999
$ws.Range("N34").Value = -6548.571599999999
$ws.Range("H39").Value = 2024.5714
$ws.Range("J39").Value = 2400.8
$ws.Range("L39").Value = 7202.400000000001
$ws.Range("N39").Value = -7790.400000000001
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H129").Value = 2231.5334
$ws.Range("I129").Value = 1249.2858
$ws.Range("J129").Value = 3091
$ws.Range("K129").Value = 3747.8574
$ws.Range("L129").Value = 9273
$ws.Range("M129").Value = 1252.1426
$ws.Range("N129").Value = -19273
$ws.Range("H131").Value = 5318.125
$ws.Range("J131").Value = 4860
$ws.Range("L131").Value = 14580
$ws.Range("N131").Value = -24660
$ws.Range("H132").Value = 2886.111
$ws.Range("J132").Value = 3757.6
$ws.Range("L132").Value = 33818.4
$ws.Range("N132").Value = -38878.4
$ws.Range("N93").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1629.5
$ws.Range("I2").Value = 261
$ws.Range("J2").Value = 2998
$ws.Range("K2").Value = 261
$ws.Range("L2").Value = 2998
$ws.Range("M2").Value = -148
$ws.Range("N2").Value = -3224
$ws.Range("H113").Value = 2741.6667
$ws.Range("J113").Value = 4926.909
$ws.Range("L113").Value = 4926.909
$ws.Range("N113").Value = -9266.909
$ws.Range("H141").Value = 63259.8
$ws.Range("J141").Value = 63259.8
$ws.Range("L141").Value = 63259.8
$ws.Range("N141").Value = -73619.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4099.9
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 4999.8335
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 4999.8335
$ws.Range("M46").Value = -2562
$ws.Range("N46").Value = -5375.8335
$ws.Range("H55").Value = 468.35715
$ws.Range("I55").Value = 158.07692
$ws.Range("J55").Value = 4502
$ws.Range("K55").Value = 158.07692
$ws.Range("L55").Value = 4502
$ws.Range("M55").Value = 14.92308
$ws.Range("N55").Value = -4848
$ws.Range("H74").Value = 43000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 43000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1246.4445
$ws.Range("I96").Value = 1280.5714
$ws.Range("J96").Value = 1127
$ws.Range("K96").Value = 1280.5714
$ws.Range("L96").Value = 1127
$ws.Range("M96").Value = 92.42859999999996
$ws.Range("N96").Value = -3873
$ws.Range("H100").Value = 2080
$ws.Range("I100").Value = 2580
$ws.Range("K100").Value = 5160
$ws.Range("M100").Value = -4619
$ws.Range("H122").Value = 607.6316
$ws.Range("I122").Value = 561.7646999999999
$ws.Range("K122").Value = 1685.2941
$ws.Range("M122").Value = 764.7059000000002
